{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2023-03-14 Tuesday\", \"2023-03-15 Wednesday\"],\n  [\"54-31=\", \"31+6=\"],\n  [\"99+0=\", \"72-55=\"],\n  [\"33-29=\", \"76-24=\"],\n  [\"23+20=\", \"42-30=\"],\n  [\"44+20=\", \"66-15=\"],\n  [\"51+43=\", \"89-4=\"],\n  [\"16-16=\", \"57-42=\"],\n  [\"63-40=\", \"57+10=\"],\n  [\"98-41=\", \"29-6=\"],\n  [\"45-16=\", \"78-66=\"],\n  [\"33-27=\", \"5+12=\"],\n  [\"77-34=\", \"23+27=\"],\n  [\"30-15=\", \"40+8=\"],\n  [\"15+59=\", \"77-43=\"],\n  [\"4+84=\", \"57-48=\"],\n  [\"67+19=\", \"13+50=\"],\n  [\"52+27=\", \"68-67=\"],\n  [\"46+36=\", \"4+2=\"],\n  [\"36+37=\", \"90-81=\"],\n  [\"34+11=\", \"31+48=\"],\n  [\"7+24=\", \"50-25=\"],\n  [\"9+90=\", \"40+7=\"],\n  [\"23+38=\", \"45-2=\"],\n  [\"31+53=\", \"4+60=\"],\n  [\"35+56=\", \"34+43=\"],\n  [\"60+21=\", \"22+14=\"],\n  [\"82+12=\", \"12+34=\"],\n  [\"68+21=\", \"39+28=\"],\n  [\"84-41=\", \"2-1=\"],\n  [\"21+63=\", \"96-15=\"],\n  [\"73+21=\", \"77-9=\"],\n  [\"76+8=\", \"66-24=\"],\n  [\"36+39=\", \"91-18=\"],\n  [\"92-64=\", \"53-38=\"],\n  [\"91+6=\", \"0+80=\"],\n  [\"66-55=\", \"40+9=\"],\n  [\"8+72=\", \"68-58=\"],\n  [\"97-94=\", \"42-21=\"],\n  [\"70-40=\", \"46-28=\"],\n  [\"60-25=\", \"69-49=\"],\n  [\"17+5=\", \"35+28=\"],\n  [\"42-28=\", \"28+60=\"],\n  [\"18+7=\", \"32-15=\"],\n  [\"79-52=\", \"87-55=\"],\n  [\"52-26=\", \"2+13=\"],\n  [\"87-1=\", \"98-1=\"],\n  [\"59-25=\", \"21-14=\"],\n  [\"64+11=\", \"19-5=\"],\n  [\"52+40=\", \"73-44=\"],\n  [\"58-51=\", \"54-31=\"],\n  [\"84+9=\", \"63-23=\"],\n  [\"18-5=\", \"99-92=\"],\n  [\"9+69=\", \"57+12=\"],\n  [\"61+9=\", \"77+16=\"],\n  [\"21+27=\", \"9-4=\"],\n  [\"6+11=\", \"8+13=\"],\n  [\"87-77=\", \"82-9=\"],\n  [\"75+11=\", \"46+27=\"],\n  [\"10+86=\", \"28+62=\"],\n  [\"50-44=\", \"32-31=\"],\n  [\"77-53=\", \"79-61=\"],\n  [\"80-61=\", \"70+8=\"],\n  [\"48-13=\", \"0+1=\"],\n  [\"97-72=\", \"92+3=\"],\n  [\"78+5=\", \"79-11=\"],\n  [\"46-45=\", \"92-88=\"],\n  [\"38+55=\", \"9+29=\"],\n  [\"56-35=\", \"58+3=\"],\n  [\"44+25=\", \"29-6=\"],\n  [\"16+77=\", \"34-32=\"],\n  [\"6+72=\", \"43-24=\"],\n  [\"9+48=\", \"53-35=\"],\n  [\"35+5=\", \"75-4=\"],\n  [\"10+87=\", \"1+3=\"],\n  [\"1+93=\", \"19+46=\"],\n  [\"57-54=\", \"18-14=\"],\n  [\"0+41=\", \"79-35=\"],\n  [\"67+3=\", \"11+37=\"],\n  [\"4+16=\", \"21-20=\"],\n  [\"84-50=\", \"29+61=\"],\n  [\"81-37=\", \"22+74=\"],\n  [\"9+66=\", \"37-16=\"],\n  [\"26-2=\", \"35+17=\"],\n  [\"49-25=\", \"17-7=\"],\n  [\"59-48=\", \"2+44=\"],\n  [\"24+71=\", \"84+12=\"],\n  [\"31-1=\", \"79-2=\"],\n  [\"99-31=\", \"2+44=\"],\n  [\"31-31=\", \"32-14=\"],\n  [\"36+12=\", \"59-12=\"],\n  [\"35+22=\", \"74-39=\"],\n  [\"20-4=\", \"61+24=\"],\n  [\"30+63=\", \"96-69=\"],\n  [\"2+94=\", \"18-0=\"],\n  [\"43+44=\", \"27-18=\"],\n  [\"16+58=\", \"61-31=\"],\n  [\"4+63=\", \"73-72=\"],\n  [\"19+73=\", \"1+1=\"],\n  [\"25-13=\", \"56+16=\"],\n  [\"55-31=\", \"80-47=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\nfunction ReplaceText($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute([Type]::Missing, $false, $false, $false, $false, $false, $true, 1, $false, [Type]::Missing, 2) | Out-Null\n}\n\nReplaceText '2023-03-14 Tuesday' '2023-03-15 Wednesday'\nReplaceText '54-31=' '31+6='\nReplaceText '99+0=' '72-55='\nReplaceText '33-29=' '76-24='\nReplaceText '23+20=' '42-30='\nReplaceText '44+20=' '66-15='\nReplaceText '51+43=' '89-4='\nReplaceText '16-16=' '57-42='\nReplaceText '63-40=' '57+10='\nReplaceText '98-41=' '29-6='\nReplaceText '45-16=' '78-66='\nReplaceText '33-27=' '5+12='\nReplaceText '77-34=' '23+27='\nReplaceText '30-15=' '40+8='\nReplaceText '15+59=' '77-43='\nReplaceText '4+84=' '57-48='\nReplaceText '67+19=' '13+50='\nReplaceText '52+27=' '68-67='\nReplaceText '46+36=' '4+2='\nReplaceText '36+37=' '90-81='\nReplaceText '34+11=' '31+48='\nReplaceText '7+24=' '50-25='\nReplaceText '9+90=' '40+7='\nReplaceText '23+38=' '45-2='\nReplaceText '31+53=' '4+60='\nReplaceText '35+56=' '34+43='\nReplaceText '60+21=' '22+14='\nReplaceText '82+12=' '12+34='\nReplaceText '68+21=' '39+28='\nReplaceText '84-41=' '2-1='\nReplaceText '21+63=' '96-15='\nReplaceText '73+21=' '77-9='\nReplaceText '76+8=' '66-24='\nReplaceText '36+39=' '91-18='\nReplaceText '92-64=' '53-38='\nReplaceText '91+6=' '0+80='\nReplaceText '66-55=' '40+9='\nReplaceText '8+72=' '68-58='\nReplaceText '97-94=' '42-21='\nReplaceText '70-40=' '46-28='\nReplaceText '60-25=' '69-49='\nReplaceText '17+5=' '35+28='\nReplaceText '42-28=' '28+60='\nReplaceText '18+7=' '32-15='\nReplaceText '79-52=' '87-55='\nReplaceText '52-26=' '2+13='\nReplaceText '87-1=' '98-1='\nReplaceText '59-25=' '21-14='\nReplaceText '64+11=' '19-5='\nReplaceText '52+40=' '73-44='\nReplaceText '58-51=' '54-31='\nReplaceText '84+9=' '63-23='\nReplaceText '18-5=' '99-92='\nReplaceText '9+69=' '57+12='\nReplaceText '61+9=' '77+16='\nReplaceText '21+27=' '9-4='\nReplaceText '6+11=' '8+13='\nReplaceText '87-77=' '82-9='\nReplaceText '75+11=' '46+27='\nReplaceText '10+86=' '28+62='\nReplaceText '50-44=' '32-31='\nReplaceText '77-53=' '79-61='\nReplaceText '80-61=' '70+8='\nReplaceText '48-13=' '0+1='\nReplaceText '97-72=' '92+3='\nReplaceText '78+5=' '79-11='\nReplaceText '46-45=' '92-88='\nReplaceText '38+55=' '9+29='\nReplaceText '56-35=' '58+3='\nReplaceText '44+25=' '29-6='\nReplaceText '16+77=' '34-32='\nReplaceText '6+72=' '43-24='\nReplaceText '9+48=' '53-35='\nReplaceText '35+5=' '75-4='\nReplaceText '10+87=' '1+3='\nReplaceText '1+93=' '19+46='\nReplaceText '57-54=' '18-14='\nReplaceText '0+41=' '79-35='\nReplaceText '67+3=' '11+37='\nReplaceText '4+16=' '21-20='\nReplaceText '84-50=' '29+61='\nReplaceText '81-37=' '22+74='\nReplaceText '9+66=' '37-16='\nReplaceText '26-2=' '35+17='\nReplaceText '49-25=' '17-7='\nReplaceText '59-48=' '2+44='\nReplaceText '24+71=' '84+12='\nReplaceText '31-1=' '79-2='\nReplaceText '99-31=' '2+44='\nReplaceText '31-31=' '32-14='\nReplaceText '36+12=' '59-12='\nReplaceText '35+22=' '74-39='\nReplaceText '20-4=' '61+24='\nReplaceText '30+63=' '96-69='\nReplaceText '2+94=' '18-0='\nReplaceText '43+44=' '27-18='\nReplaceText '16+58=' '61-31='\nReplaceText '4+63=' '73-72='\nReplaceText '19+73=' '1+1='\nReplaceText '25-13=' '56+16='\nReplaceText '55-31=' '80-47='\n"}
